# "Make sure Ben's update is in other folder."
#
# The canonical-OOXML diff for this commit touches exactly eight
# <a14:hiddenLine> extension blocks living inside the freeform "Freeform 6"
# decorative shape on eight slide layouts (slideLayout1/2/4/5/6/8/13/14.xml).
# In every single hunk the *only* change is that the two namespace
# declarations on that one element get swapped:
#
#   <a14:hiddenLine xmlns:a14="...2010/main" xmlns="" w="9525">
#   -> <a14:hiddenLine xmlns="" xmlns:a14="...2010/main" w="9525">
#
# Nothing inside the element (fill color, cap style, arrowheads, the
# w="9525" value, ...) changes, and no other byte in the package changes.
# That <a14:hiddenLine> blob is PowerPoint's cache of a shape's outline
# formatting while Shape.Line.Visible is False; it is written once by
# whatever tool last touched that shape and is otherwise carried through
# untouched. This PowerPoint session never toggled those shapes' outlines
# (Ben's edit was unrelated "move the deck into the other folder" work),
# so there is no Shape/Line/Fill/etc. object-model mutation that would
# legitimately cause PowerPoint to rewrite that already-hidden line cache
# here - touching Line/Fill only rewrites the live <a:ln> sibling, never
# the preserved <a14:hiddenLine> extension (verified empirically: setting
# Line.Visible/Weight/ForeColor, Fill, Shadow, Glow, rotation, duplicating
# or recreating the shape, resaving the whole deck, etc. all leave that
# extension's bytes, including attribute order, untouched).
#
# So there is nothing for this session to legitimately change: the eight
# shapes' outline formatting is not being edited, just the file's location
# in the repo. Leave the presentation exactly as-is.
